$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "CP_lista_campos_registro"
$ws.Range("E6").Value = "Como se aprecia en la imagen X, la lista no trae informacion y al ser un campo obligatorio no puese realizar la accion "
$ws.Range("G6").Value = "abierto"
$ws.Range("I6").Value = "alta"
$ws.Range("J6").Value = "sebastian"

$ws.Range("L16").Select()
